$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 118
$ws.Range("D118").Value = 44511
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 9000
$ws.Range("L118").Value = 9500
$ws.Range("M118").Value = 9250
$ws.Range("O118").Value = 'Región de Ñuble'
$ws.Range("P118").Value = 462

# Row 119
$ws.Range("D119").Value = 44511
$ws.Range("I119").Value = 'Segunda'
$ws.Range("J119").Value = 500
$ws.Range("K119").Value = 8000
$ws.Range("L119").Value = 8000
$ws.Range("M119").Value = 8000
$ws.Range("O119").Value = 'Región de Ñuble'
$ws.Range("P119").Value = 400

# Row 120
$ws.Range("D120").Value = 44306
$ws.Range("I120").Value = 'Primera'
$ws.Range("J120").Value = 400
$ws.Range("K120").Value = 6500
$ws.Range("L120").Value = 6500
$ws.Range("M120").Value = 6500
$ws.Range("O120").Value = 'Región de Ñuble'
$ws.Range("P120").Value = 325

# Row 121
$ws.Range("D121").Value = 44306
$ws.Range("I121").Value = 'Segunda'
$ws.Range("J121").Value = 400
$ws.Range("K121").Value = 5000
$ws.Range("L121").Value = 5000
$ws.Range("M121").Value = 5000
$ws.Range("O121").Value = 'Región de Ñuble'
$ws.Range("P121").Value = 250

# Row 122
$ws.Range("D122").Value = 44211
$ws.Range("I122").Value = 'Primera'
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 6500
$ws.Range("M122").Value = 6250
$ws.Range("O122").Value = 'Región de Ñuble'
$ws.Range("P122").Value = 312

# Row 123
$ws.Range("D123").Value = 44211
$ws.Range("I123").Value = 'Segunda'
$ws.Range("J123").Value = 300
$ws.Range("K123").Value = 5000
$ws.Range("L123").Value = 5000
$ws.Range("M123").Value = 5000
$ws.Range("O123").Value = 'Región de Ñuble'
$ws.Range("P123").Value = 250

# Row 124
$ws.Range("D124").Value = 44215
$ws.Range("I124").Value = 'Primera'
$ws.Range("J124").Value = 400
$ws.Range("K124").Value = 6000
$ws.Range("L124").Value = 6000
$ws.Range("M124").Value = 6000
$ws.Range("O124").Value = 'Región de Ñuble'
$ws.Range("P124").Value = 300

# Row 125
$ws.Range("D125").Value = 44215
$ws.Range("I125").Value = 'Segunda'
$ws.Range("J125").Value = 400
$ws.Range("K125").Value = 5000
$ws.Range("L125").Value = 5000
$ws.Range("M125").Value = 5000
$ws.Range("O125").Value = 'Región de Ñuble'
$ws.Range("P125").Value = 250

# Row 126
$ws.Range("D126").Value = 44504
$ws.Range("I126").Value = 'Primera'
$ws.Range("J126").Value = 350
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7000
$ws.Range("M126").Value = 6571
$ws.Range("O126").Value = 'Chillán'
$ws.Range("P126").Value = 329

# Row 127
$ws.Range("D127").Value = 44505
$ws.Range("I127").Value = 'Primera'
$ws.Range("J127").Value = 250
$ws.Range("K127").Value = 8000
$ws.Range("L127").Value = 8500
$ws.Range("M127").Value = 8200
$ws.Range("O127").Value = 'Provincia del Elquí'
$ws.Range("P127").Value = 410

# Row 128
$ws.Range("D128").Value = 44425
$ws.Range("I128").Value = 'Primera'
$ws.Range("J128").Value = 600
$ws.Range("K128").Value = 5000
$ws.Range("L128").Value = 5500
$ws.Range("M128").Value = 5250
$ws.Range("O128").Value = 'Región de Ñuble'
$ws.Range("P128").Value = 262

# Row 129
$ws.Range("D129").Value = 44425
$ws.Range("I129").Value = 'Segunda'
$ws.Range("J129").Value = 300
$ws.Range("K129").Value = 4500
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 4500
$ws.Range("O129").Value = 'Región de Ñuble'
$ws.Range("P129").Value = 225

# Row 130
$ws.Range("D130").Value = 44343
$ws.Range("I130").Value = 'Primera'
$ws.Range("J130").Value = 1000
$ws.Range("K130").Value = 4500
$ws.Range("L130").Value = 5000
$ws.Range("M130").Value = 4750
$ws.Range("O130").Value = 'Región de Ñuble'
$ws.Range("P130").Value = 238

# Row 131
$ws.Range("D131").Value = 44343
$ws.Range("I131").Value = 'Segunda'
$ws.Range("J131").Value = 500
$ws.Range("K131").Value = 4000
$ws.Range("L131").Value = 4000
$ws.Range("M131").Value = 4000
$ws.Range("O131").Value = 'Región de Ñuble'
$ws.Range("P131").Value = 200

# Row 132
$ws.Range("D132").Value = 44370
$ws.Range("I132").Value = 'Primera'
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 5000
$ws.Range("M132").Value = 4750
$ws.Range("O132").Value = 'Región de Ñuble'
$ws.Range("P132").Value = 238

# Row 133
$ws.Range("D133").Value = 44392
$ws.Range("I133").Value = 'Primera'
$ws.Range("J133").Value = 600
$ws.Range("K133").Value = 5000
$ws.Range("L133").Value = 5500
$ws.Range("M133").Value = 5250
$ws.Range("O133").Value = 'Región de Ñuble'
$ws.Range("P133").Value = 262

# Row 134
$ws.Range("D134").Value = 44392
$ws.Range("I134").Value = 'Segunda'
$ws.Range("J134").Value = 300
$ws.Range("K134").Value = 4000
$ws.Range("L134").Value = 4000
$ws.Range("M134").Value = 4000
$ws.Range("O134").Value = 'Región de Ñuble'
$ws.Range("P134").Value = 200

# Row 135
$ws.Range("D135").Value = 44295
$ws.Range("I135").Value = 'Primera'
$ws.Range("J135").Value = 600
$ws.Range("K135").Value = 5500
$ws.Range("L135").Value = 6000
$ws.Range("M135").Value = 5750
$ws.Range("O135").Value = 'Región de Ñuble'
$ws.Range("P135").Value = 288

# Row 136
$ws.Range("D136").Value = 44295
$ws.Range("I136").Value = 'Segunda'
$ws.Range("J136").Value = 300
$ws.Range("K136").Value = 5000
$ws.Range("L136").Value = 5000
$ws.Range("M136").Value = 5000
$ws.Range("O136").Value = 'Región de Ñuble'
$ws.Range("P136").Value = 250

# Row 137
$ws.Range("D137").Value = 44217
$ws.Range("I137").Value = 'Primera'
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 6500
$ws.Range("L137").Value = 7000
$ws.Range("M137").Value = 6750
$ws.Range("O137").Value = 'Región de Ñuble'
$ws.Range("P137").Value = 338

# Row 138
$ws.Range("D138").Value = 44217
$ws.Range("I138").Value = 'Segunda'
$ws.Range("J138").Value = 500
$ws.Range("K138").Value = 5500
$ws.Range("L138").Value = 5500
$ws.Range("M138").Value = 5500
$ws.Range("O138").Value = 'Región de Ñuble'
$ws.Range("P138").Value = 275

# Row 139
$ws.Range("D139").Value = 44509
$ws.Range("I139").Value = 'Primera'
$ws.Range("J139").Value = 800
$ws.Range("K139").Value = 6500
$ws.Range("L139").Value = 7000
$ws.Range("M139").Value = 6750
$ws.Range("O139").Value = 'Región Metropolitana'
$ws.Range("P139").Value = 338

# Row 140
$ws.Range("D140").Value = 44509
$ws.Range("I140").Value = 'Segunda'
$ws.Range("J140").Value = 400
$ws.Range("K140").Value = 5500
$ws.Range("L140").Value = 5500
$ws.Range("M140").Value = 5500
$ws.Range("O140").Value = 'Región Metropolitana'
$ws.Range("P140").Value = 275

# Row 141
$ws.Range("D141").Value = 44421
$ws.Range("I141").Value = 'Primera'
$ws.Range("J141").Value = 500
$ws.Range("K141").Value = 5000
$ws.Range("L141").Value = 5500
$ws.Range("M141").Value = 5200
$ws.Range("O141").Value = 'Región de Ñuble'
$ws.Range("P141").Value = 260

# Row 142
$ws.Range("D142").Value = 44421
$ws.Range("I142").Value = 'Segunda'
$ws.Range("J142").Value = 300
$ws.Range("K142").Value = 4500
$ws.Range("L142").Value = 4500
$ws.Range("M142").Value = 4500
$ws.Range("O142").Value = 'Región de Ñuble'
$ws.Range("P142").Value = 225

# Row 143
$ws.Range("D143").Value = 44383
$ws.Range("I143").Value = 'Primera'
$ws.Range("J143").Value = 600
$ws.Range("K143").Value = 5000
$ws.Range("L143").Value = 5500
$ws.Range("M143").Value = 5250
$ws.Range("O143").Value = 'Región de Ñuble'
$ws.Range("P143").Value = 262

# Row 144
$ws.Range("D144").Value = 44383
$ws.Range("I144").Value = 'Segunda'
$ws.Range("J144").Value = 300
$ws.Range("K144").Value = 4500
$ws.Range("L144").Value = 4500
$ws.Range("M144").Value = 4500
$ws.Range("O144").Value = 'Región de Ñuble'
$ws.Range("P144").Value = 225

# Row 145
$ws.Range("D145").Value = 44307
$ws.Range("I145").Value = 'Primera'
$ws.Range("J145").Value = 200
$ws.Range("K145").Value = 6000
$ws.Range("L145").Value = 6500
$ws.Range("M145").Value = 6250
$ws.Range("O145").Value = 'Región de Ñuble'
$ws.Range("P145").Value = 312

# Row 146
$ws.Range("D146").Value = 44307
$ws.Range("I146").Value = 'Segunda'
$ws.Range("J146").Value = 100
$ws.Range("K146").Value = 5000
$ws.Range("L146").Value = 5000
$ws.Range("M146").Value = 5000
$ws.Range("O146").Value = 'Región de Ñuble'
$ws.Range("P146").Value = 250

# Row 147
$ws.Range("A147").Value = 11
$ws.Range("B147").Value = 'Vega Monumental Concepción'
$ws.Range("C147").Value = 'Bíobío'
$ws.Range("D147").Value = 44433
$ws.Range("E147").Value = 8
$ws.Range("F147").Value = 100114013
$ws.Range("G147").Value = 'Zanahoria'
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 600
$ws.Range("K147").Value = 5000
$ws.Range("L147").Value = 5500
$ws.Range("M147").Value = 5250
$ws.Range("N147").Value = '$/saco 20 kilos'
$ws.Range("O147").Value = 'Región de Ñuble'
$ws.Range("P147").Value = 262
$ws.Range("Q147").Value = 20
$ws.Range("R147").Value = 'Hortaliza'

# Row 148
$ws.Range("A148").Value = 11
$ws.Range("B148").Value = 'Vega Monumental Concepción'
$ws.Range("C148").Value = 'Bíobío'
$ws.Range("D148").Value = 44433
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = 100114013
$ws.Range("G148").Value = 'Zanahoria'
$ws.Range("H148").Value = 'Sin especificar'
$ws.Range("I148").Value = 'Segunda'
$ws.Range("J148").Value = 300
$ws.Range("K148").Value = 4500
$ws.Range("L148").Value = 4500
$ws.Range("M148").Value = 4500
$ws.Range("N148").Value = '$/saco 20 kilos'
$ws.Range("O148").Value = 'Región de Ñuble'
$ws.Range("P148").Value = 225
$ws.Range("Q148").Value = 20
$ws.Range("R148").Value = 'Hortaliza'

